$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 76
$ws.Range("B4").Value = 71
$ws.Range("B5").Value = 64
$ws.Range("B6").Value = 71
$ws.Range("B7").Value = 69
$ws.Range("B9").Value = 72
$ws.Range("B11").Value = 78
$ws.Range("B13").Value = 71
$ws.Range("B14").Value = 72
